$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear everything currently used (A1:J2) so stale cells (row 2, columns D-J) go away
$ws.Range("A1:J2").Clear()

# Write the new header row
$ws.Range("A1").Value = "Pregunta"
$ws.Range("B1").Value = "Respuesta"
$ws.Range("C1").Value = "Conteo"
